$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 0.2341817432078095
$ws.Range("D2").Value = 0.06982524292032366
$ws.Range("E2").Value = 0.1253305760472827
$ws.Range("F2").Value = 1.962084723588688
$ws.Range("G2").Value = 0.002520600672214817
$ws.Range("J2").Value = 0.2151374932702907
$ws.Range("L2").Value = 0.2092131224034972
$ws.Range("M2").Value = 11.28904589160032
$ws.Range("O2").Value = 5.627354544310833

$ws.Range("C3").Value = 0.2436480373668957
$ws.Range("D3").Value = 0.06995959182046363
$ws.Range("E3").Value = 0.1235244002138245
$ws.Range("F3").Value = 2.02086215390591
$ws.Range("G3").Value = 0.002526810755308015
$ws.Range("J3").Value = 0.2096047403423853
$ws.Range("L3").Value = 0.1930001562005259
$ws.Range("M3").Value = 9.986315243197339
$ws.Range("O3").Value = 5.779349527788185

$ws.Range("C4").Value = 0.2497722846182704
$ws.Range("D4").Value = 0.07007460212251004
$ws.Range("E4").Value = 0.1225436138629767
$ws.Range("F4").Value = 2.059480990279475
$ws.Range("G4").Value = 0.002530789394027733
$ws.Range("J4").Value = 0.206470090037314
$ws.Range("L4").Value = 0.1831296581791406
$ws.Range("M4").Value = 9.183018286896015
$ws.Range("O4").Value = 5.880022600035943

$ws.Range("C5").Value = 0.2523455570065032
$ws.Range("D5").Value = 0.07012963562143426
$ws.Range("E5").Value = 0.1221756374686578
$ws.Range("F5").Value = 2.075847419913941
$ws.Range("G5").Value = 0.002532452585709733
$ws.Range("J5").Value = 0.2052572485916215
$ws.Range("L5").Value = 0.1791279315229986
$ws.Range("M5").Value = 8.854784017118902
$ws.Range("O5").Value = 5.922874022840375

$ws.Range("C6").Value = 0.2527775089704889
$ws.Range("D6").Value = 0.07013926681808158
$ws.Range("E6").Value = 0.122116434871046
$ws.Range("F6").Value = 2.078602834025695
$ws.Range("G6").Value = 0.002532731292257036
$ws.Range("J6").Value = 0.2050597152387894
$ws.Range("L6").Value = 0.1784646741243563
$ws.Range("M6").Value = 8.800226895146977
$ws.Range("O6").Value = 5.930099155832778

$ws.Range("C7").Value = 0.2498066756007766
$ws.Range("D7").Value = 0.07007531127307232
$ws.Range("E7").Value = 0.1225385235313787
$ws.Range("F7").Value = 2.059699176391291
$ws.Range("G7").Value = 0.002530811654446604
$ws.Range("J7").Value = 0.206453473712898
$ws.Range("L7").Value = 0.1830756069309416
$ws.Range("M7").Value = 9.178595215779353
$ws.Range("O7").Value = 5.880593142506143

$ws.Range("C8").Value = 0.2373804788559717
$ws.Range("D8").Value = 0.06986480898724778
$ws.Range("E8").Value = 0.1246808612610337
$ws.Range("F8").Value = 1.981821976646749
$ws.Range("G8").Value = 0.002522707667086954
$ws.Range("J8").Value = 0.2131744227470733
$ws.Range("L8").Value = 0.2036049858354829
$ws.Range("M8").Value = 10.84055618897429
$ws.Range("O8").Value = 5.678224602089301

$ws.Range("C9").Value = 0.2155251940569141
$ws.Range("D9").Value = 0.06971079326163476
$ws.Range("E9").Value = 0.1299255793787282
$ws.Range("F9").Value = 1.849490129974157
$ws.Range("G9").Value = 0.002508119615227428
$ws.Range("J9").Value = 0.2285071445250821
$ws.Range("L9").Value = 0.2445662340460757
$ws.Range("M9").Value = 14.07407331500991
$ws.Range("O9").Value = 5.340636090594643

$ws.Range("C10").Value = 0.2010517124399556
$ws.Range("D10").Value = 0.06975678624282722
$ws.Range("E10").Value = 0.13445451043156
$ws.Range("F10").Value = 1.765129866053464
$ws.Range("G10").Value = 0.002498181643311892
$ws.Range("J10").Value = 0.2411893051781675
$ws.Range("L10").Value = 0.2751451279691253
$ws.Range("M10").Value = 16.43677071472075
$ws.Range("O10").Value = 5.130056073368337

$ws.Range("C11").Value = 0.1948227474371009
$ws.Range("D11").Value = 0.06981265317483576
$ws.Range("E11").Value = 0.1366709002605049
$ws.Range("F11").Value = 1.729641499581554
$ws.Range("G11").Value = 0.002493826576420355
$ws.Range("J11").Value = 0.2472911638120081
$ws.Range("L11").Value = 0.2891755437177324
$ws.Range("M11").Value = 17.50948335486021
$ws.Range("O11").Value = 5.042679883802322

$ws.Range("C12").Value = 0.1925160683321572
$ws.Range("D12").Value = 0.06983886942894202
$ws.Range("E12").Value = 0.1375334610003449
$ws.Range("F12").Value = 1.716626417709563
$ws.Range("G12").Value = 0.002492200990835972
$ws.Range("J12").Value = 0.2496517651566137
$ws.Range("L12").Value = 0.2945069545765193
$ws.Range("M12").Value = 17.91544843165747
$ws.Range("O12").Value = 5.010828053318164

$ws.Range("C13").Value = 0.1930105192743063
$ws.Range("D13").Value = 0.06983299759998829
$ws.Range("E13").Value = 0.1373466451699983
$ws.Range("F13").Value = 1.71941047394413
$ws.Range("G13").Value = 0.002492550045002391
$ws.Range("J13").Value = 0.2491411103993784
$ws.Range("L13").Value = 0.2933579022468109
$ws.Range("M13").Value = 17.82802672074649
$ws.Range("O13").Value = 5.017632553024043

$ws.Range("C14").Value = 0.1946319254903166
$ws.Range("D14").Value = 0.06981470835952308
$ws.Range("E14").Value = 0.1367413926565959
$ws.Range("F14").Value = 1.728562203001715
$ws.Range("G14").Value = 0.002493692367316177
$ws.Range("J14").Value = 0.247484358066373
$ws.Range("L14").Value = 0.2896137864455852
$ws.Range("M14").Value = 17.54288699438058
$ws.Range("O14").Value = 5.040034504477973

$ws.Range("C15").Value = 0.1956318990140993
$ws.Range("D15").Value = 0.06980416585098936
$ws.Range("E15").Value = 0.1363737121742261
$ws.Range("F15").Value = 1.73422331685228
$ws.Range("G15").Value = 0.002494395136323452
$ws.Range("J15").Value = 0.2464761206459087
$ws.Range("L15").Value = 0.2873228391090379
$ws.Range("M15").Value = 17.36820004442518
$ws.Range("O15").Value = 5.053918019187591

$ws.Range("C16").Value = 0.2014660229906404
$ws.Range("D16").Value = 0.06975384133114915
$ws.Range("E16").Value = 0.1343128743091029
$ws.Range("F16").Value = 1.76750793026298
$ws.Range("G16").Value = 0.002498469571328155
$ws.Range("J16").Value = 0.24079740735489
$ws.Range("L16").Value = 0.2742307167068816
$ws.Range("M16").Value = 16.36662874237373
$ws.Range("O16").Value = 5.135937839203535

$ws.Range("C17").Value = 0.205136743204676
$ws.Range("D17").Value = 0.06973194300866936
$ws.Range("E17").Value = 0.1330891716736105
$ws.Range("F17").Value = 1.78867242993234
$ws.Range("G17").Value = 0.002501011384381304
$ws.Range("J17").Value = 0.2374003799916125
$ws.Range("L17").Value = 0.2662306276127708
$ws.Range("M17").Value = 15.75169765115129
$ws.Range("O17").Value = 5.188428071975011

$ws.Range("C18").Value = 0.2072814041530151
$ws.Range("D18").Value = 0.06972263441530657
$ws.Range("E18").Value = 0.1323999773288875
$ws.Range("F18").Value = 1.801116968903258
$ws.Range("G18").Value = 0.002502488986093094
$ws.Range("J18").Value = 0.2354776791243864
$ws.Range("L18").Value = 0.2616404007798394
$ws.Range("M18").Value = 15.397805336614
$ws.Range("O18").Value = 5.219409844777488

$ws.Range("C19").Value = 0.2080132480617287
$ws.Range("D19").Value = 0.06972004616876148
$ws.Range("E19").Value = 0.1321691197923016
$ws.Range("F19").Value = 1.805376816181493
$ws.Range("G19").Value = 0.002502991966785972
$ws.Range("J19").Value = 0.2348319780381445
$ws.Range("L19").Value = 0.260088121159697
$ws.Range("M19").Value = 15.27794757279986
$ws.Range("O19").Value = 5.230034799083825

$ws.Range("C20").Value = 0.2047425287719058
$ws.Range("D20").Value = 0.06973393371348635
$ws.Range("E20").Value = 0.1332179153001647
$ws.Range("F20").Value = 1.78639129646529
$ws.Range("G20").Value = 0.002500739189468562
$ws.Range("J20").Value = 0.2377587574138147
$ws.Range("L20").Value = 0.2670810822374108
$ws.Range("M20").Value = 15.81717847214856
$ws.Range("O20").Value = 5.182758373750175

$ws.Range("C21").Value = 0.1941542569702435
$ws.Range("D21").Value = 0.06981994270110903
$ws.Range("E21").Value = 0.1369185319984183
$ws.Range("F21").Value = 1.725862554580871
$ws.Range("G21").Value = 0.002493356201486827
$ws.Range("J21").Value = 0.2479696129944671
$ws.Range("L21").Value = 0.2907130150547061
$ws.Range("M21").Value = 17.62664572688431
$ws.Range("O21").Value = 5.033420770134569

$ws.Range("C22").Value = 0.1875383572841827
$ws.Range("D22").Value = 0.06990567165048134
$ws.Range("E22").Value = 0.1394730964881248
$ws.Range("F22").Value = 1.688776555869453
$ws.Range("G22").Value = 0.002488668334045311
$ws.Range("J22").Value = 0.2549352716944639
$ws.Range("L22").Value = 0.3062656862172446
$ws.Range("M22").Value = 18.80782028963392
$ws.Range("O22").Value = 4.943034322178534

$ws.Range("C23").Value = 0.1910412092043785
$ws.Range("D23").Value = 0.06985720274399654
$ws.Range("E23").Value = 0.1380969544841406
$ws.Range("F23").Value = 1.708340956981679
$ws.Range("G23").Value = 0.002491157856039355
$ws.Range("J23").Value = 0.2511900935841282
$ws.Range("L23").Value = 0.2979546593233806
$ws.Range("M23").Value = 18.17751608191048
$ws.Range("O23").Value = 4.990606700124374

$ws.Range("C24").Value = 0.204920646311578
$ws.Range("D24").Value = 0.06973302349867083
$ws.Range("E24").Value = 0.1331596657366489
$ws.Range("F24").Value = 1.787421735409907
$ws.Range("G24").Value = 0.002500862198013822
$ws.Range("J24").Value = 0.2375966407769567
$ws.Range("L24").Value = 0.2666965633464429
$ws.Range("M24").Value = 15.78757571895494
$ws.Range("O24").Value = 5.18531914057013

$ws.Range("C25").Value = 0.2211633632879355
$ws.Range("D25").Value = 0.06972462529002144
$ws.Range("E25").Value = 0.1283910977651104
$ws.Range("F25").Value = 1.88306236855346
$ws.Range("G25").Value = 0.002511928008835214
$ws.Range("J25").Value = 0.2241183078025983
$ws.Range("L25").Value = 0.2334043952599387
$ws.Range("M25").Value = 13.20180630487567
$ws.Range("O25").Value = 5.425483725511555

